# Insert a new bad-word entry ("cho'choq" / id 368 / severity 50)
# right before the existing "chumo" row (row 104), shifting all
# subsequent rows down by one (old row 104 -> new row 105, etc.).
# This grows the sheet from A1:C195 to A1:C196.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 104; everything currently at
# row 104 and below shifts down to row 105 and below.
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row with the new entry.
$ws.Cells.Item(104, 1).Value = 368
$ws.Cells.Item(104, 2).Value = "cho'choq"
$ws.Cells.Item(104, 3).Value = 50
